$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 114 ----
$ws.Cells.Item(113, 1).Copy()
$ws.Cells.Item(114, 1).PasteSpecial(-4122)   # xlPasteFormats, reuse the existing date style
$excel.CutCopyMode = $false

$ws.Cells.Item(114, 1).Value = 45448.2916666667
$ws.Cells.Item(114, 2).Value = 0
$ws.Cells.Item(114, 3).Value = 3.24000000953674
$ws.Cells.Item(114, 4).Value = 3.24000000953674
$ws.Cells.Item(114, 5).Value = 3.24000000953674
$ws.Cells.Item(114, 6).Value = 3.24000000953674

$ws.Cells.Item(114, 7).NumberFormat = "@"
$ws.Cells.Item(114, 7).Value = "3.24000000953674"
$ws.Cells.Item(114, 7).ClearFormats()

$ws.Cells.Item(114, 8).Value = "AGAIN.MI"

# ---- Row 115 ----
$ws.Cells.Item(113, 1).Copy()
$ws.Cells.Item(115, 1).PasteSpecial(-4122)   # xlPasteFormats, reuse the existing date style
$excel.CutCopyMode = $false

$ws.Cells.Item(115, 1).Value = 45449.5224189815
$ws.Cells.Item(115, 2).Value = 2000
$ws.Cells.Item(115, 3).Value = 3.29999995231628
$ws.Cells.Item(115, 4).Value = 3.22000002861023
$ws.Cells.Item(115, 5).Value = 3.29999995231628
$ws.Cells.Item(115, 6).Value = 3.22000002861023

$ws.Cells.Item(115, 7).NumberFormat = "@"
$ws.Cells.Item(115, 7).Value = "3.22000002861023"
$ws.Cells.Item(115, 7).ClearFormats()

$ws.Cells.Item(115, 8).Value = "AGAIN.MI"
